$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.289.20"
$ws.Range("E2").Value = "  +4.02%  "

$ws.Range("D3").Value = "1.713.34"
$ws.Range("E3").Value = "  +3.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.15"
$ws.Range("E5").Value = "  +4.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5409"
$ws.Range("E6").Value = "  +3.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2761"
$ws.Range("E8").Value = "  +3.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06493"
$ws.Range("E9").Value = "  +2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.61"
$ws.Range("E10").Value = "  +4.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.736"
$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("D13").Value = "1.713.55"
$ws.Range("E13").Value = "  +3.41%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.951.89"
$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6058"
$ws.Range("E15").Value = "  +7.09%  "

$ws.Range("D16").Value = "0.0₅8307"
$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.01"
$ws.Range("E17").Value = "  +5.48%  "

$ws.Range("D18").Value = "27.253.09"
$ws.Range("E18").Value = "  +3.89%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.803"
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.33"
$ws.Range("E20").Value = "  +9.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.97"
$ws.Range("E22").Value = "  +6.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.221"
$ws.Range("E23").Value = "  +3.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.38"
$ws.Range("E25").Value = "  +2.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1261"
$ws.Range("E26").Value = "  +4.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.442"
$ws.Range("E27").Value = "  +2.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.62"
$ws.Range("E28").Value = "  +3.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.603"
$ws.Range("E29").Value = "  +6.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05632"
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.324"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.641"
$ws.Range("E32").Value = "  +3.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.511"
$ws.Range("E33").Value = "  +3.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.641"
$ws.Range("E34").Value = "  +3.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9778"
$ws.Range("E35").Value = "  +3.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.835"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.422"
$ws.Range("E37").Value = "  +0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5837"
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01644"
$ws.Range("E39").Value = "  +2.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.908"
$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("D41").Value = "1.051.04"
$ws.Range("E41").Value = "  +1.41%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9999"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8366"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.62"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").Value = "1.856.83"
$ws.Range("E45").Value = "  +3.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.93"
$ws.Range("E46").Value = "  +2.48%  "

$ws.Range("D47").Value = "0.0₈109"
$ws.Range("E47").Value = "  +4.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.095"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9907"
$ws.Range("E49").Value = "  -1.23%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4342"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05242"
$ws.Range("E51").Value = "  -1.41%  "
